$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("A2").Value = "ECs"
    $ws.Range("B2").Value = "Vtn"
    $ws.Range("C2").Value = "Tnfrsf11b"
    $ws.Range("D2").Value = "FAPs"
    $ws.Range("E2").Value = 3
    $ws.Range("F2").Value = 1
    $ws.Range("G2").Value = 6.201571333333333
    $ws.Range("H2").Value = 18.604714
    $ws.Range("I2").Value = 0.05221490529364391
    $ws.Range("J2").Value = 0.07406232529850043
    $ws.Range("K2").Value = 3
    $ws.Range("L2").Value = 1
    $ws.Range("M2").Value = 2.311298
    $ws.Range("N2").Value = 6.933894
    $ws.Range("O2").Value = 0.8122763614007964
    $ws.Range("P2").Value = 0.863868705896424
    $ws.Range("Q2").Value = 14.33367941959067
    $ws.Range("R2").Value = 129.003114776316
    $ws.Range("S2").Value = 0.04241293328280826
    $ws.Range("T2").Value = 0.06398012511129554

    # Row 3
    $ws.Range("A3").Value = "ECs"
    $ws.Range("B3").Value = "Vtn"
    $ws.Range("C3").Value = "Tnfrsf11b"
    $ws.Range("D3").Value = "M1"
    $ws.Range("E3").Value = 3
    $ws.Range("F3").Value = 1
    $ws.Range("G3").Value = 6.201571333333333
    $ws.Range("H3").Value = 18.604714
    $ws.Range("I3").Value = 0.05221490529364391
    $ws.Range("J3").Value = 0.07406232529850043
    $ws.Range("K3").Value = 1
    $ws.Range("L3").Value = 0.3333333333333333
    $ws.Range("M3").Value = 0.02434666666666667
    $ws.Range("N3").Value = 0.07304
    $ws.Range("O3").Value = 0.008556327142686946
    $ws.Range("P3").Value = 0.009099788701510982
    $ws.Range("Q3").Value = 0.1509875900622222
    $ws.Range("R3").Value = 1.35888831056
    $ws.Range("S3").Value = 0.0004467678114168337
    $ws.Range("T3").Value = 0.0006739515109589252

    # Row 4
    $ws.Range("A4").Value = "ECs"
    $ws.Range("B4").Value = "Vtn"
    $ws.Range("C4").Value = "Tnfrsf11b"
    $ws.Range("D4").Value = "sCs"
    $ws.Range("E4").Value = 3
    $ws.Range("F4").Value = 1
    $ws.Range("G4").Value = 6.201571333333333
    $ws.Range("H4").Value = 18.604714
    $ws.Range("I4").Value = 0.05221490529364391
    $ws.Range("J4").Value = 0.07406232529850043
    $ws.Range("K4").Value = 2
    $ws.Range("L4").Value = 1
    $ws.Range("M4").Value = 0.509813
    $ws.Range("N4").Value = 1.019626
    $ws.Range("O4").Value = 0.1791673114565167
    $ws.Range("P4").Value = 0.1270315054020651
    $ws.Range("Q4").Value = 3.161641686160666
    $ws.Range("R4").Value = 18.969850116964
    $ws.Range("S4").Value = 0.009355204199418819
    $ws.Range("T4").Value = 0.009408248676245963

    # Row 5
    $ws.Range("A5").Value = "FAPs"
    $ws.Range("B5").Value = "Vtn"
    $ws.Range("C5").Value = "Tnfrsf11b"
    $ws.Range("D5").Value = "FAPs"
    $ws.Range("E5").Value = 3
    $ws.Range("F5").Value = 1
    $ws.Range("G5").Value = 7.461641333333333
    $ws.Range("H5").Value = 22.384924
    $ws.Range("I5").Value = 0.06282422221945559
    $ws.Range("J5").Value = 0.0891107233935555
    $ws.Range("K5").Value = 3
    $ws.Range("L5").Value = 1
    $ws.Range("M5").Value = 2.311298
    $ws.Range("N5").Value = 6.933894
    $ws.Range("O5").Value = 0.8122763614007964
    $ws.Range("P5").Value = 0.863868705896424
    $ws.Range("Q5").Value = 17.24607669045067
    $ws.Range("R5").Value = 155.214690214056
    $ws.Range("S5").Value = 0.05103063063225445
    $ws.Range("T5").Value = 0.07697996529948498

    # Row 6
    $ws.Range("A6").Value = "FAPs"
    $ws.Range("B6").Value = "Vtn"
    $ws.Range("C6").Value = "Tnfrsf11b"
    $ws.Range("D6").Value = "M1"
    $ws.Range("E6").Value = 3
    $ws.Range("F6").Value = 1
    $ws.Range("G6").Value = 7.461641333333333
    $ws.Range("H6").Value = 22.384924
    $ws.Range("I6").Value = 0.06282422221945559
    $ws.Range("J6").Value = 0.0891107233935555
    $ws.Range("K6").Value = 1
    $ws.Range("L6").Value = 0.3333333333333333
    $ws.Range("M6").Value = 0.02434666666666667
    $ws.Range("N6").Value = 0.07304
    $ws.Range("O6").Value = 0.008556327142686946
    $ws.Range("P6").Value = 0.009099788701510982
    $ws.Range("Q6").Value = 0.1816660943288889
    $ws.Range("R6").Value = 1.63499484896
    $ws.Range("S6").Value = 0.0005375445977945242
    $ws.Range("T6").Value = 0.0008108887539201466

    # Row 7
    $ws.Range("A7").Value = "FAPs"
    $ws.Range("B7").Value = "Vtn"
    $ws.Range("C7").Value = "Tnfrsf11b"
    $ws.Range("D7").Value = "sCs"
    $ws.Range("E7").Value = 3
    $ws.Range("F7").Value = 1
    $ws.Range("G7").Value = 7.461641333333333
    $ws.Range("H7").Value = 22.384924
    $ws.Range("I7").Value = 0.06282422221945559
    $ws.Range("J7").Value = 0.0891107233935555
    $ws.Range("K7").Value = 2
    $ws.Range("L7").Value = 1
    $ws.Range("M7").Value = 0.509813
    $ws.Range("N7").Value = 1.019626
    $ws.Range("O7").Value = 0.1791673114565167
    $ws.Range("P7").Value = 0.1270315054020651
    $ws.Range("Q7").Value = 3.804041753070666
    $ws.Range("R7").Value = 22.82425051842399
    $ws.Range("S7").Value = 0.01125604698940662
    $ws.Range("T7").Value = 0.01131986934015038

    # Row 8
    $ws.Range("A8").Value = "sCs"
    $ws.Range("B8").Value = "Vtn"
    $ws.Range("C8").Value = "Tnfrsf11b"
    $ws.Range("D8").Value = "FAPs"
    $ws.Range("E8").Value = 2
    $ws.Range("F8").Value = 1
    $ws.Range("G8").Value = 105.106922
    $ws.Range("H8").Value = 210.213844
    $ws.Range("I8").Value = 0.8849608724869005
    $ws.Range("J8").Value = 0.836826951307944
    $ws.Range("K8").Value = 3
    $ws.Range("L8").Value = 1
    $ws.Range("M8").Value = 2.311298
    $ws.Range("N8").Value = 6.933894
    $ws.Range("O8").Value = 0.8122763614007964
    $ws.Range("P8").Value = 0.863868705896424
    $ws.Range("Q8").Value = 242.933418604756
    $ws.Range("R8").Value = 1457.600511628536
    $ws.Range("S8").Value = 0.7188327974857337
    $ws.Range("T8").Value = 0.7229086154856434

    # Row 9
    $ws.Range("A9").Value = "sCs"
    $ws.Range("B9").Value = "Vtn"
    $ws.Range("C9").Value = "Tnfrsf11b"
    $ws.Range("D9").Value = "M1"
    $ws.Range("E9").Value = 2
    $ws.Range("F9").Value = 1
    $ws.Range("G9").Value = 105.106922
    $ws.Range("H9").Value = 210.213844
    $ws.Range("I9").Value = 0.8849608724869005
    $ws.Range("J9").Value = 0.836826951307944
    $ws.Range("K9").Value = 1
    $ws.Range("L9").Value = 0.3333333333333333
    $ws.Range("M9").Value = 0.02434666666666667
    $ws.Range("N9").Value = 0.07304
    $ws.Range("O9").Value = 0.008556327142686946
    $ws.Range("P9").Value = 0.009099788701510982
    $ws.Range("Q9").Value = 2.559003194293333
    $ws.Range("R9").Value = 15.35401916576
    $ws.Range("S9").Value = 0.007572014733475589
    $ws.Range("T9").Value = 0.00761494843663191

    # Row 10
    $ws.Range("A10").Value = "sCs"
    $ws.Range("B10").Value = "Vtn"
    $ws.Range("C10").Value = "Tnfrsf11b"
    $ws.Range("D10").Value = "sCs"
    $ws.Range("E10").Value = 2
    $ws.Range("F10").Value = 1
    $ws.Range("G10").Value = 105.106922
    $ws.Range("H10").Value = 210.213844
    $ws.Range("I10").Value = 0.8849608724869005
    $ws.Range("J10").Value = 0.836826951307944
    $ws.Range("K10").Value = 2
    $ws.Range("L10").Value = 1
    $ws.Range("M10").Value = 0.509813
    $ws.Range("N10").Value = 1.019626
    $ws.Range("O10").Value = 0.1791673114565167
    $ws.Range("P10").Value = 0.1270315054020651
    $ws.Range("Q10").Value = 53.584875225586
    $ws.Range("R10").Value = 214.339500902344
    $ws.Range("S10").Value = 0.1585560602676912
    $ws.Range("T10").Value = 0.1063033873856688

